# Update the "cryptos" price list: refresh Price/Volume(1h) figures for
# most rows, and fix the coin ordering for two swapped pairs
# (WrappedBTC/BitcoinCash at rows 16-17, USDe/PEPE at rows 47-48).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading "'" forces Excel to store the literal as text (matching the
# original inline-string cell type) instead of auto-parsing look-alikes
# such as "1.00" or "3.219.14" as numbers; resetting the style afterwards
# clears the quote-prefix formatting flag Excel would otherwise persist.
$ws.Range("D2").Value = "'66.395.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.32%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.257.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +7.20%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'582.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +5.05%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'152.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +7.80%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.250.20"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +7.27%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +5.79%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +9.83%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +6.81%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.490"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +5.43%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'37.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.97%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +6.54%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.791.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +7.33%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'WrappedBTC"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'66.638.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.49%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'BitcoinCash"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'555.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +13.36%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.262.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +6.67%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +2.88%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +6.45%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +6.31%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.746"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +8.52%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +10.45%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'13.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +7.44%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'81.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.72%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.25%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +19.27%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'3.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +9.35%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +7.25%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'27.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +7.27%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +7.03%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.43%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +6.45%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'568.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +8.68%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +4.19%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'6.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +7.30%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'55.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +5.35%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0454"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +11.92%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.0867"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +8.41%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.131"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +8.18%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +10.29%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.219.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +11.45%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.07%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +15.54%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +10.74%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'26.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +6.42%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'USDe"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.07%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'PEPE"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.0₃0559"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +4.39%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'125.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.80%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +4.45%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +9.06%  "
$ws.Range("E51").Style = "Normal"
